# One-way flight scenario is added
# The "Stays" sheet test data is updated so several test cases point at a
# one-way-flight destination city (Belgrade / Zagreb) instead of the old
# ski-resort locations (Kopaonik / Kikinda). The "Stays" sheet also becomes
# the active/selected sheet (instead of "Flights"), with a new active
# selection cell.

$wb = $excel.ActiveWorkbook
$flights = $wb.Worksheets.Item("Flights")
$stays = $wb.Worksheets.Item("Stays")

# Update the "location" column (B) test data for TC_001..TC_004
$stays.Range("B4").Value = "Belgrade"
$stays.Range("B5").Value = "Zagreb"
$stays.Range("B6").Value = "Belgrade"
$stays.Range("B7").Value = "Belgrade"

# Update the matching "verifResults" column (K) test data
$stays.Range("K4").Value = "Belgrade"
$stays.Range("K5").Value = "Zagreb"

# Make "Stays" the active sheet and update its active selection cell
$stays.Activate()
$stays.Range("L10").Select()
